# Fruta / hortaliza, semanal
#
# The sheet holds weekly price observations for "Betarraga" at Feria
# Lagunitas de Puerto Montt. A new week's observation is inserted at the
# top of the data block (row 108); every subsequent observation (rows
# 108..174) shifts down by one row, and the observation that used to be
# in the last row (174) is appended as a brand-new last row (175).
#
# Columns A, B, C, E, F, G, H, I, R are constant across this entire block
# (same market/category/quality/classification for every row), so only
# D (Fecha), J (Volumen), K (Precio minimo), L (Precio maximo),
# M (Precio promedio ponderado), N (Unidad de comercializacion),
# O (Origen), P (Precio $/Kg) and Q (Kg o Unidades) actually need to move.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$shiftCols = @("D","J","K","L","M","N","O","P","Q")
$firstRow = 108
$lastRow = 174
$newLastRow = 175

# Snapshot every value that will move, before writing anything, since the
# shift would otherwise clobber data we still need to read.
$old = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowVals = @{}
    foreach ($c in $shiftCols) {
        $rowVals[$c] = $ws.Range("$c$r").Value2()
    }
    $old[$r] = $rowVals
}

# New observation for the top row of the block; this week's data is brand
# new, everything else on row 108 (unit, origin, per-kilo price, etc.)
# stays exactly as it already was.
$ws.Range("D$firstRow").Value2 = 44438
$ws.Range("J$firstRow").Value2 = 500

# Every other existing row takes on the values that used to belong to the
# row above it (push everything down by one row).
for ($r = $firstRow + 1; $r -le $lastRow; $r++) {
    $src = $old[$r - 1]
    foreach ($c in $shiftCols) {
        $ws.Range("$c$r").Value2 = $src[$c]
    }
}

# Append a brand-new last row carrying what used to be the final
# observation (old row 174), including the columns that are constant
# across the whole block.
$ws.Range("A$newLastRow").Value2 = $ws.Range("A$lastRow").Value2()
$ws.Range("B$newLastRow").Value2 = $ws.Range("B$lastRow").Value2()
$ws.Range("C$newLastRow").Value2 = $ws.Range("C$lastRow").Value2()
$ws.Range("E$newLastRow").Value2 = $ws.Range("E$lastRow").Value2()
$ws.Range("F$newLastRow").Value2 = $ws.Range("F$lastRow").Value2()
$ws.Range("G$newLastRow").Value2 = $ws.Range("G$lastRow").Value2()
$ws.Range("H$newLastRow").Value2 = $ws.Range("H$lastRow").Value2()
$ws.Range("I$newLastRow").Value2 = $ws.Range("I$lastRow").Value2()
$ws.Range("R$newLastRow").Value2 = $ws.Range("R$lastRow").Value2()

$lastOld = $old[$lastRow]
foreach ($c in $shiftCols) {
    $ws.Range("$c$newLastRow").Value2 = $lastOld[$c]
}

# Match the date cell's formatting to the rest of the D column so it
# reuses the existing date style instead of acquiring a new one.
$ws.Range("D$newLastRow").NumberFormat = $ws.Range("D$lastRow").NumberFormat()
